# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").Value = "'76.559.95"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -0.13%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.944.28"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.01%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'198.89"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.11%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'594.77"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  +0.04%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.549"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -1.14%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.196"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.71%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'2.944.98"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +1.71%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.442"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +11.14%  "

$ws.Range("E12").Value = "  +0.31%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'3.494.90"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +1.97%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'4.88"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -0.89%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'76.603.03"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.01%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'28.27"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +2.67%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'0.0000188"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.59%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'2.964.23"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +2.36%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'13.52"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +7.34%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'8.67"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -3.93%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'373.89"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.47%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +4.00%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'2.25"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.96%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'72.32"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'3.109.56"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +2.18%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'4.27"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +1.02%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'9.61"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.45%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'0.0000107"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +1.73%  "

$ws.Range("E30").Value = "  +0.06%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'8.29"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +6.77%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.37"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -3.57%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'498.64"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -2.45%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'1.83"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.78%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'164.96"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.112"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +20.92%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'20.15"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.393"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +13.56%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'19.94"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("E41").Value = "  -6.38%  "

$ws.Range("E42").Value = "  +0.01%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'180.25"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.01%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'4.91"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -3.89%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'1.64"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -2.61%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'40.14"
$ws.Range("D46").Style = $origStyle

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'1.18"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -4.32%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.589"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.27%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'3.88"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +2.97%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = "'2.30"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -2.85%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'22.36"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +3.14%  "
